$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions - copy formatting (style) from H1 so I1/J1 match the
# existing bold/centered/bordered header style, then set their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows additions for columns I (I0) and J (IF)
$data = @(
    @(7, 8),
    @(4, 5),
    @(7, 8),
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(6, 7),
    @(8, 9),
    @(7, 7),
    @(6, 7),
    @(5, 6),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
